$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1238.5
$ws.Range("I19").Value = 1152.25
$ws.Range("J19").Value = 1324.75
$ws.Range("K19").Value = 1152.25
$ws.Range("L19").Value = 1324.75
$ws.Range("M19").Value = -977.25
$ws.Range("N19").Value = -1674.75
$ws.Range("H55").Value = 923.5
$ws.Range("I55").Value = 818.2
$ws.Range("J55").Value = 998.7143
$ws.Range("K55").Value = 818.2
$ws.Range("L55").Value = 998.7143
$ws.Range("M55").Value = -604.2
$ws.Range("N55").Value = -1426.7143
$ws.Range("H103").Value = 2138.1538
$ws.Range("J103").Value = 1805.2222
$ws.Range("L103").Value = 5415.6666
$ws.Range("N103").Value = -6587.6666
$ws.Range("H109").Value = 99597.5
$ws.Range("J109").Value = 99597.5
$ws.Range("L109").Value = 99597.5
$ws.Range("N109").Value = -102371.5
$ws.Range("H138").Value = 2506.2454
$ws.Range("I138").Value = 2053.7715
$ws.Range("J138").Value = 3386.0557
$ws.Range("K138").Value = 6161.314499999999
$ws.Range("L138").Value = 10158.1671
$ws.Range("M138").Value = -1021.314499999999
$ws.Range("N138").Value = -20438.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2083.0625
$ws.Range("I110").Value = 1805.8334
$ws.Range("J110").Value = 2914.75
$ws.Range("K110").Value = 1805.8334
$ws.Range("L110").Value = 2914.75
$ws.Range("M110").Value = 239.1666
$ws.Range("N110").Value = -7004.75
$ws.Range("H132").Value = 45079.707
$ws.Range("I132").Value = 45079.707
$ws.Range("K132").Value = 135239.121
$ws.Range("M132").Value = -132709.121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 125947
$ws.Range("I22").Value = 125947
$ws.Range("K22").Value = 125947
$ws.Range("M22").Value = -125774
$ws.Range("H107").Value = 1808.6666
$ws.Range("I107").Value = 1146
$ws.Range("K107").Value = 1146
$ws.Range("M107").Value = 774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1800.7333
$ws.Range("I16").Value = 1800.1538
$ws.Range("K16").Value = 1800.1538
$ws.Range("M16").Value = -1513.1538
$ws.Range("H31").Value = 4063.8
$ws.Range("I31").Value = 2930.6667
$ws.Range("J31").Value = 4701.1875
$ws.Range("K31").Value = 2930.6667
$ws.Range("L31").Value = 4701.1875
$ws.Range("M31").Value = -2635.6667
$ws.Range("N31").Value = -5291.1875
$ws.Range("H34").Value = 4063.8
$ws.Range("I34").Value = 2930.6667
$ws.Range("J34").Value = 4701.1875
$ws.Range("K34").Value = 2930.6667
$ws.Range("L34").Value = 4701.1875
$ws.Range("M34").Value = -2728.6667
$ws.Range("N34").Value = -5105.1875
$ws.Range("H37").Value = 2449.4
$ws.Range("J37").Value = 2449.4
$ws.Range("L37").Value = 2449.4
$ws.Range("N37").Value = -2663.4
$ws.Range("H113").Value = 1800.7333
$ws.Range("I113").Value = 1800.1538
$ws.Range("K113").Value = 1800.1538
$ws.Range("M113").Value = 369.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 404.75
$ws.Range("I97").Value = 226.33333
$ws.Range("J97").Value = 583.1667
$ws.Range("K97").Value = 678.99999
$ws.Range("L97").Value = 1749.5001
$ws.Range("M97").Value = -182.99999
$ws.Range("N97").Value = -2741.5001
$ws.Range("H98").Value = 1117.3334
$ws.Range("I98").Value = 765.6
$ws.Range("J98").Value = 1557
$ws.Range("K98").Value = 2296.8
$ws.Range("L98").Value = 4671
$ws.Range("M98").Value = -798.8000000000002
$ws.Range("N98").Value = -7667
$ws.Range("H109").Value = 4217.4287
$ws.Range("J109").Value = 9500
$ws.Range("L109").Value = 28500
$ws.Range("N109").Value = -30580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H102").Value = 3500.9048
$ws.Range("I102").Value = 2553.8667
$ws.Range("K102").Value = 2553.8667
$ws.Range("M102").Value = -931.8667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3923.6667
$ws.Range("I7").Value = 2708.4
$ws.Range("K7").Value = 2708.4
$ws.Range("M7").Value = -2596.4
$ws.Range("H40").Value = 7002.778
$ws.Range("I40").Value = 7667.6665
$ws.Range("J40").Value = 6670.3335
$ws.Range("K40").Value = 7667.6665
$ws.Range("L40").Value = 6670.3335
$ws.Range("M40").Value = -7531.6665
$ws.Range("N40").Value = -6942.3335
$ws.Range("H45").Value = 40999.8
$ws.Range("I45").Value = 29999.5
$ws.Range("K45").Value = 29999.5
$ws.Range("M45").Value = -29592.5
$ws.Range("H100").Value = 3257.4167
$ws.Range("J100").Value = 3033
$ws.Range("L100").Value = 3033
$ws.Range("N100").Value = -4115
$ws.Range("H112").Value = 29387
$ws.Range("J112").Value = 29387
$ws.Range("L112").Value = 29387
$ws.Range("N112").Value = -32341
$ws.Range("H122").Value = 4462.4736
$ws.Range("I122").Value = 3970.5557
$ws.Range("J122").Value = 4905.2
$ws.Range("K122").Value = 11911.6671
$ws.Range("L122").Value = 14715.6
$ws.Range("M122").Value = -9461.667099999999
$ws.Range("N122").Value = -19615.6
$ws.Range("H126").Value = 3923.6667
$ws.Range("I126").Value = 2708.4
$ws.Range("K126").Value = 8125.200000000001
$ws.Range("M126").Value = -5655.200000000001
$ws.Range("H132").Value = 40972.625
$ws.Range("I132").Value = 52245
$ws.Range("K132").Value = 156735
$ws.Range("M132").Value = -154205
$ws.Range("H134").Value = 79992.5
$ws.Range("J134").Value = 79992.5
$ws.Range("L134").Value = 79992.5
$ws.Range("N134").Value = -90132.5
$ws.Range("H135").Value = 51473.332
$ws.Range("J135").Value = 51473.332
$ws.Range("L135").Value = 51473.332
$ws.Range("N135").Value = -61613.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14366
$ws.Range("J41").Value = 14366
$ws.Range("L41").Value = 14366
$ws.Range("N41").Value = -15146
$ws.Range("H100").Value = 987.06665
$ws.Range("I100").Value = 817.3333
$ws.Range("K100").Value = 1634.6666
$ws.Range("M100").Value = -1093.6666
$ws.Range("H107").Value = 1488.5883
$ws.Range("I107").Value = 1077.6
$ws.Range("J107").Value = 2075.7144
$ws.Range("K107").Value = 3232.8
$ws.Range("L107").Value = 6227.1432
$ws.Range("M107").Value = -1312.8
$ws.Range("N107").Value = -10067.1432
$ws.Range("H109").Value = 96749
$ws.Range("J109").Value = 96749
$ws.Range("L109").Value = 96749
$ws.Range("N109").Value = -99523
